$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Motiviator"
$ws.Range("B3").Value = "MOTIVIATOR"
$ws.Range("C3").Value = "'123456789"
$ws.Range("D3").Value = "terreneitor"
$ws.Range("E3").Value = "motiviator@motiviator.com"
$ws.Range("F3").Value = 2
